# Generate Report for Handoff
#
# Updates the localization-status workbook to reflect a fresh handoff report:
#  - Status moves from "Handed back: in sync with en-US" to "In Translation"
#    on the Overview sheet (zh-cn + de-de columns) and on each language
#    sheet's Status column.
#  - The various "Latest ... Datetime" timestamps advance a couple of
#    minutes.
#  - Each language sheet's "Error Detail" cell now reports that the handback
#    file version is stale, with a link to the current vs. latest commit.
#  - A couple of columns are resized to fit their new content.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2cf57d1b06786932c84a63bc7ea85d8af122cb88/e2e/9e653d48-c409-4f43-b574-36ba8ffd9fd1.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cda1e62299376a0ebdc4630d5c5529b2dab15488/e2e/9e653d48-c409-4f43-b574-36ba8ffd9fd1.md."

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-11-15 16:18:19"

$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-11-15 16:18:04"
$zhcn.Range("P2").Value = $errorDetail

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-11-15 16:18:19"
$dede.Range("P2").Value = $errorDetail

$dede.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
